# Update Leve price/profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets to match refreshed market-board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 26: Everything Is Impossible / Budding Ash Wand
$ws.Range("H26").Value = 34606.25
$ws.Range("J26").Value = 34606.25
$ws.Range("L26").Value = 34606.25
$ws.Range("N26").Value = -35294.25
# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Range("H28").Value = 511.90323
$ws.Range("I28").Value = 626.625
$ws.Range("J28").Value = 118.57143
$ws.Range("K28").Value = 626.625
$ws.Range("L28").Value = 118.57143
$ws.Range("M28").Value = -141.625
$ws.Range("N28").Value = -1088.57143
# Row 58: A Matter of Vital Importance / Mega-Potion of Vitality
$ws.Range("H58").Value = 1149.2307
$ws.Range("I58").Value = 490.375
$ws.Range("J58").Value = 2203.4
$ws.Range("K58").Value = 1471.125
$ws.Range("L58").Value = 6610.200000000001
$ws.Range("M58").Value = -1321.125
$ws.Range("N58").Value = -6910.200000000001
# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 2375.8235
$ws.Range("I100").Value = 2866.6667
$ws.Range("J100").Value = 2108.0908
$ws.Range("K100").Value = 2866.6667
$ws.Range("L100").Value = 2108.0908
$ws.Range("M100").Value = -2325.6667
$ws.Range("N100").Value = -3190.0908
# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1039.1515
$ws.Range("J112").Value = 1039.1515
$ws.Range("L112").Value = 3117.4545
$ws.Range("N112").Value = -5333.4545
# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 2612.125
$ws.Range("I113").Value = 2632.3333
$ws.Range("J113").Value = 2600
$ws.Range("K113").Value = 2632.3333
$ws.Range("L113").Value = 2600
$ws.Range("M113").Value = 621.6667000000002
$ws.Range("N113").Value = -9108
# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 377688.25
$ws.Range("I129").Value = 1092234.2
$ws.Range("J129").Value = 3402.238
$ws.Range("K129").Value = 3276702.6
$ws.Range("L129").Value = 10206.714
$ws.Range("M129").Value = -3271702.6
$ws.Range("N129").Value = -20206.714
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 20231.326
$ws.Range("I132").Value = 3875.3076
$ws.Range("J132").Value = 69299.38
$ws.Range("K132").Value = 11625.9228
$ws.Range("L132").Value = 207898.14
$ws.Range("M132").Value = -9095.9228
$ws.Range("N132").Value = -212958.14
# Row 133: Big Brush, Big Dreams / Ginseng Angle Brush
$ws.Range("H133").Value = 73410.71000000001
$ws.Range("J133").Value = 73410.71000000001
$ws.Range("L133").Value = 73410.71000000001
$ws.Range("N133").Value = -83530.71000000001
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2522.803
$ws.Range("J138").Value = 2638.535
$ws.Range("L138").Value = 7915.605
$ws.Range("N138").Value = -18195.605
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 2845.4546
$ws.Range("I2").Value = 2830
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 2830
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -2717
$ws.Range("N2").Value = -3226
# Row 4: Eyes Bigger than the Plate / Bronze Plate
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 34569.96
$ws.Range("I32").Value = 35374.48
$ws.Range("J32").Value = 14457
$ws.Range("K32").Value = 35374.48
$ws.Range("L32").Value = 14457
$ws.Range("M32").Value = -35087.48
$ws.Range("N32").Value = -15031
# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 2845.4546
$ws.Range("I116").Value = 2830
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2830
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -536
$ws.Range("N116").Value = -7588
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 2845.4546
$ws.Range("I3").Value = 2830
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 2830
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -2716
$ws.Range("N3").Value = -3228
# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 735.2593000000001
$ws.Range("I94").Value = 663.2632
$ws.Range("K94").Value = 663.2632
$ws.Range("M94").Value = -212.2632
# Row 114: Halfhearted Effort / Bluespirit Halfheart Saw
$ws.Range("H114").Value = 19999.889
$ws.Range("J114").Value = 19999.889
$ws.Range("L114").Value = 19999.889
$ws.Range("N114").Value = -28677.889
# Row 117: Idol Hands / Titanbronze Chakrams
$ws.Range("H117").Value = 49540.4
$ws.Range("J117").Value = 49540.4
$ws.Range("L117").Value = 49540.4
$ws.Range("N117").Value = -58718.4
# Row 130: Annals of the Empire I / Chondrite Magitek Axe
$ws.Range("H130").Value = 48372.8
$ws.Range("J130").Value = 48372.8
$ws.Range("L130").Value = 48372.8
$ws.Range("N130").Value = -58412.8
# Row 141: Awl Dreams Come True / Ra'Kaznar Awl
$ws.Range("H141").Value = 41916.6
$ws.Range("J141").Value = 41916.6
$ws.Range("L141").Value = 41916.6
$ws.Range("N141").Value = -52276.6
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 4584.586
$ws.Range("I31").Value = 1778.2059
$ws.Range("J31").Value = 7235.0557
$ws.Range("K31").Value = 1778.2059
$ws.Range("L31").Value = 7235.0557
$ws.Range("M31").Value = -1483.2059
$ws.Range("N31").Value = -7825.0557
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 4584.586
$ws.Range("I34").Value = 1778.2059
$ws.Range("J34").Value = 7235.0557
$ws.Range("K34").Value = 1778.2059
$ws.Range("L34").Value = 7235.0557
$ws.Range("M34").Value = -1576.2059
$ws.Range("N34").Value = -7639.0557
# Row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 1910.3636
$ws.Range("I86").Value = 1922.8
$ws.Range("J86").Value = 1900
$ws.Range("K86").Value = 1922.8
$ws.Range("L86").Value = 1900
$ws.Range("M86").Value = -799.8
$ws.Range("N86").Value = -4146
# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 1910.3636
$ws.Range("I89").Value = 1922.8
$ws.Range("J89").Value = 1900
$ws.Range("K89").Value = 9614
$ws.Range("L89").Value = 9500
$ws.Range("M89").Value = -3998
$ws.Range("N89").Value = -20732
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 50880.863
$ws.Range("I132").Value = 2067.5881
$ws.Range("K132").Value = 6202.7643
$ws.Range("M132").Value = -3672.7643
$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water / Boiled Egg
$ws.Range("H4").Value = 2106
$ws.Range("I4").Value = 188.5
$ws.Range("J4").Value = 3008.353
$ws.Range("K4").Value = 565.5
$ws.Range("L4").Value = 9025.059000000001
$ws.Range("M4").Value = -453.5
$ws.Range("N4").Value = -9249.059000000001
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 1242.8182
$ws.Range("I5").Value = 722.46155
$ws.Range("J5").Value = 1994.4445
$ws.Range("K5").Value = 2167.38465
$ws.Range("L5").Value = 5983.333500000001
$ws.Range("M5").Value = -2055.38465
$ws.Range("N5").Value = -6207.333500000001
# Row 11: Putting the Squeeze On / Orange Juice
$ws.Range("H11").Value = 2595.532
$ws.Range("I11").Value = 2837.8572
$ws.Range("J11").Value = 560
$ws.Range("K11").Value = 8513.571599999999
$ws.Range("L11").Value = 1680
$ws.Range("M11").Value = -8373.571599999999
$ws.Range("N11").Value = -1960
# Row 98: Sweet Kiss of Death / Rice Vinegar
$ws.Range("H98").Value = 2234.6667
$ws.Range("I98").Value = 2000
$ws.Range("J98").Value = 2352
$ws.Range("K98").Value = 6000
$ws.Range("L98").Value = 7056
$ws.Range("M98").Value = -4502
$ws.Range("N98").Value = -10052
# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 3549.7354
$ws.Range("I113").Value = 4388.346
$ws.Range("J113").Value = 824.25
$ws.Range("K113").Value = 13165.038
$ws.Range("L113").Value = 2472.75
$ws.Range("M113").Value = -10995.038
$ws.Range("N113").Value = -6812.75
# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 3262.5122
$ws.Range("I122").Value = 723.2162
$ws.Range("J122").Value = 26751
$ws.Range("K122").Value = 6508.9458
$ws.Range("L122").Value = 240759
$ws.Range("M122").Value = -4058.9458
$ws.Range("N122").Value = -245659
# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 1242.8182
$ws.Range("I135").Value = 722.46155
$ws.Range("J135").Value = 1994.4445
$ws.Range("K135").Value = 6502.15395
$ws.Range("L135").Value = 17950.0005
$ws.Range("M135").Value = -3967.15395
$ws.Range("N135").Value = -23020.0005
# Row 140: Sweet, Sweet Bean Juice / Mesquite Juice
$ws.Range("H140").Value = 2352.2
$ws.Range("I140").Value = 1690.25
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 5070.75
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = 109.25
$ws.Range("N140").Value = -25360
$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me / Bone Hora
$ws.Range("H5").Value = 15000
$ws.Range("J5").Value = 15500
$ws.Range("L5").Value = 15500
$ws.Range("N5").Value = -15724
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2224.7805
$ws.Range("I132").Value = 1756.2413
$ws.Range("K132").Value = 5268.7239
$ws.Range("M132").Value = -2738.7239
# Row 140: The Right Rod / Ra'Kaznar Rod
$ws.Range("H140").Value = 36499.75
$ws.Range("J140").Value = 36499.75
$ws.Range("L140").Value = 36499.75
$ws.Range("N140").Value = -46859.75
$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 19120
$ws.Range("I136").Value = 53271.43
$ws.Range("J136").Value = 2441.3953
$ws.Range("K136").Value = 159814.29
$ws.Range("L136").Value = 7324.1859
$ws.Range("M136").Value = -157264.29
$ws.Range("N136").Value = -12424.1859
